$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(3, 8).Value = 40000
$ws.Cells.Item(3, 10).Value = 40000
$ws.Cells.Item(3, 12).Value = 40000
$ws.Cells.Item(3, 14).Value = -40228
$ws.Cells.Item(8, 8).Value = 202.33333
$ws.Cells.Item(8, 10).Value = 99
$ws.Cells.Item(8, 12).Value = 297
$ws.Cells.Item(8, 14).Value = -575
$ws.Cells.Item(9, 8).Value = 5000000
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 13).Value = $null
$ws.Cells.Item(17, 8).Value = 665.13635
$ws.Cells.Item(17, 10).Value = 665.13635
$ws.Cells.Item(17, 12).Value = 1995.40905
$ws.Cells.Item(17, 14).Value = -2331.40905
$ws.Cells.Item(40, 8).Value = 2893.2
$ws.Cells.Item(40, 9).Value = 2847.4285
$ws.Cells.Item(40, 10).Value = 3000
$ws.Cells.Item(40, 11).Value = 2847.4285
$ws.Cells.Item(40, 12).Value = 3000
$ws.Cells.Item(40, 13).Value = -2672.4285
$ws.Cells.Item(40, 14).Value = -3350
$ws.Cells.Item(53, 8).Value = 699.2857
$ws.Cells.Item(53, 9).Value = 482.5
$ws.Cells.Item(53, 10).Value = 2000
$ws.Cells.Item(53, 11).Value = 482.5
$ws.Cells.Item(53, 12).Value = 2000
$ws.Cells.Item(53, 13).Value = 154.5
$ws.Cells.Item(53, 14).Value = -3274
$ws.Cells.Item(69, 8).Value = 7399.909
$ws.Cells.Item(69, 10).Value = 8231.666999999999
$ws.Cells.Item(69, 12).Value = 24695.001
$ws.Cells.Item(69, 14).Value = -26443.001
$ws.Cells.Item(70, 8).Value = 5916.6
$ws.Cells.Item(70, 10).Value = 5916.6
$ws.Cells.Item(70, 12).Value = 17749.8
$ws.Cells.Item(70, 14).Value = -18289.8
$ws.Cells.Item(72, 8).Value = 7399.909
$ws.Cells.Item(72, 10).Value = 8231.666999999999
$ws.Cells.Item(72, 12).Value = 74085.003
$ws.Cells.Item(72, 14).Value = -82821.003
$ws.Cells.Item(73, 8).Value = 5916.6
$ws.Cells.Item(73, 10).Value = 5916.6
$ws.Cells.Item(73, 12).Value = 17749.8
$ws.Cells.Item(73, 14).Value = -19621.8
$ws.Cells.Item(76, 8).Value = 5069.364
$ws.Cells.Item(76, 9).Value = 4194.3335
$ws.Cells.Item(76, 10).Value = 5397.5
$ws.Cells.Item(76, 11).Value = 4194.3335
$ws.Cells.Item(76, 12).Value = 5397.5
$ws.Cells.Item(76, 13).Value = -3879.3335
$ws.Cells.Item(76, 14).Value = -6027.5
$ws.Cells.Item(79, 8).Value = 5069.364
$ws.Cells.Item(79, 9).Value = 4194.3335
$ws.Cells.Item(79, 10).Value = 5397.5
$ws.Cells.Item(79, 11).Value = 4194.3335
$ws.Cells.Item(79, 12).Value = 5397.5
$ws.Cells.Item(79, 13).Value = -3102.3335
$ws.Cells.Item(79, 14).Value = -7581.5
$ws.Cells.Item(86, 8).Value = 3586.889
$ws.Cells.Item(86, 9).Value = 3468.2856
$ws.Cells.Item(86, 10).Value = 4002
$ws.Cells.Item(86, 11).Value = 3468.2856
$ws.Cells.Item(86, 12).Value = 4002
$ws.Cells.Item(86, 13).Value = -2345.2856
$ws.Cells.Item(86, 14).Value = -6248
$ws.Cells.Item(87, 8).Value = 20000
$ws.Cells.Item(87, 10).Value = 20000
$ws.Cells.Item(87, 12).Value = 20000
$ws.Cells.Item(87, 14).Value = -22496
$ws.Cells.Item(89, 8).Value = 3586.889
$ws.Cells.Item(89, 9).Value = 3468.2856
$ws.Cells.Item(89, 10).Value = 4002
$ws.Cells.Item(89, 11).Value = 17341.428
$ws.Cells.Item(89, 12).Value = 20010
$ws.Cells.Item(89, 13).Value = -11725.428
$ws.Cells.Item(89, 14).Value = -31242
$ws.Cells.Item(90, 8).Value = 20000
$ws.Cells.Item(90, 10).Value = 20000
$ws.Cells.Item(90, 12).Value = 60000
$ws.Cells.Item(90, 14).Value = -72480
$ws.Cells.Item(92, 8).Value = 5115.375
$ws.Cells.Item(92, 9).Value = 1744
$ws.Cells.Item(92, 10).Value = 6239.1665
$ws.Cells.Item(92, 11).Value = 1744
$ws.Cells.Item(92, 12).Value = 6239.1665
$ws.Cells.Item(92, 13).Value = -496
$ws.Cells.Item(92, 14).Value = -8735.166499999999
$ws.Cells.Item(95, 8).Value = 46330.668
$ws.Cells.Item(95, 10).Value = 46330.668
$ws.Cells.Item(95, 12).Value = 46330.668
$ws.Cells.Item(95, 14).Value = -51822.668
$ws.Cells.Item(96, 8).Value = 967
$ws.Cells.Item(96, 9).Value = 1810
$ws.Cells.Item(96, 10).Value = 124
$ws.Cells.Item(96, 11).Value = 5430
$ws.Cells.Item(96, 12).Value = 372
$ws.Cells.Item(96, 13).Value = -4057
$ws.Cells.Item(96, 14).Value = -3118
$ws.Cells.Item(97, 8).Value = 1670
$ws.Cells.Item(97, 9).Value = 1500
$ws.Cells.Item(97, 11).Value = 4500
$ws.Cells.Item(97, 13).Value = -4004
$ws.Cells.Item(100, 8).Value = 801.2
$ws.Cells.Item(100, 9).Value = 777.5
$ws.Cells.Item(100, 11).Value = 777.5
$ws.Cells.Item(100, 13).Value = -236.5
$ws.Cells.Item(102, 8).Value = 40000
$ws.Cells.Item(102, 10).Value = 40000
$ws.Cells.Item(102, 12).Value = 40000
$ws.Cells.Item(102, 14).Value = -46490
$ws.Cells.Item(103, 8).Value = 634.7273
$ws.Cells.Item(103, 9).Value = 597.1667
$ws.Cells.Item(103, 10).Value = 679.8
$ws.Cells.Item(103, 11).Value = 1791.5001
$ws.Cells.Item(103, 12).Value = 2039.4
$ws.Cells.Item(103, 13).Value = -1205.5001
$ws.Cells.Item(103, 14).Value = -3211.4
$ws.Cells.Item(105, 8).Value = 70000
$ws.Cells.Item(105, 10).Value = 70000
$ws.Cells.Item(105, 12).Value = 70000
$ws.Cells.Item(105, 14).Value = -76988
$ws.Cells.Item(106, 8).Value = 1700
$ws.Cells.Item(106, 10).Value = 0
$ws.Cells.Item(106, 12).Value = 0
$ws.Cells.Item(106, 14).Value = $null
$ws.Cells.Item(112, 8).Value = 1531.4584
$ws.Cells.Item(112, 9).Value = 1199.7693
$ws.Cells.Item(112, 10).Value = 1923.4546
$ws.Cells.Item(112, 11).Value = 3599.3079
$ws.Cells.Item(112, 12).Value = 5770.3638
$ws.Cells.Item(112, 13).Value = -2491.3079
$ws.Cells.Item(112, 14).Value = -7986.3638
$ws.Cells.Item(132, 8).Value = 771719.5600000001
$ws.Cells.Item(132, 9).Value = 2759.4546
$ws.Cells.Item(132, 10).Value = 5001000
$ws.Cells.Item(132, 11).Value = 8278.363799999999
$ws.Cells.Item(132, 12).Value = 15003000
$ws.Cells.Item(132, 13).Value = -5748.363799999999
$ws.Cells.Item(132, 14).Value = -15008060
$ws.Cells.Item(135, 8).Value = 4414.148
$ws.Cells.Item(135, 9).Value = 687.3200000000001
$ws.Cells.Item(135, 10).Value = 50999.5
$ws.Cells.Item(135, 11).Value = 6185.88
$ws.Cells.Item(135, 12).Value = 458995.5
$ws.Cells.Item(135, 13).Value = -3650.88
$ws.Cells.Item(135, 14).Value = -464065.5
$ws.Cells.Item(137, 8).Value = 1421.6923
$ws.Cells.Item(137, 9).Value = 854.25
$ws.Cells.Item(137, 11).Value = 2562.75
$ws.Cells.Item(137, 13).Value = -12.75
$ws.Cells.Item(138, 8).Value = 2361.9663
$ws.Cells.Item(138, 9).Value = 2679.8928
$ws.Cells.Item(138, 10).Value = 2216.0327
$ws.Cells.Item(138, 11).Value = 8039.678400000001
$ws.Cells.Item(138, 12).Value = 6648.098100000001
$ws.Cells.Item(138, 13).Value = -2899.678400000001
$ws.Cells.Item(138, 14).Value = -16928.0981
$ws.Cells.Item(141, 8).Value = 5038.2666
$ws.Cells.Item(141, 9).Value = 2951.7273
$ws.Cells.Item(141, 10).Value = 10776.25
$ws.Cells.Item(141, 11).Value = 8855.1819
$ws.Cells.Item(141, 12).Value = 32328.75
$ws.Cells.Item(141, 13).Value = -3675.1819
$ws.Cells.Item(141, 14).Value = -42688.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2269
$ws.Cells.Item(2, 9).Value = 2674.75
$ws.Cells.Item(2, 10).Value = 1457.5
$ws.Cells.Item(2, 11).Value = 2674.75
$ws.Cells.Item(2, 12).Value = 1457.5
$ws.Cells.Item(2, 13).Value = -2561.75
$ws.Cells.Item(2, 14).Value = -1683.5
$ws.Cells.Item(5, 8).Value = 1017.1667
$ws.Cells.Item(5, 9).Value = 620.6
$ws.Cells.Item(5, 11).Value = 620.6
$ws.Cells.Item(5, 13).Value = -508.6
$ws.Cells.Item(26, 8).Value = 450
$ws.Cells.Item(26, 9).Value = 450
$ws.Cells.Item(26, 11).Value = 450
$ws.Cells.Item(26, 13).Value = -120
$ws.Cells.Item(32, 8).Value = 7648.4136
$ws.Cells.Item(32, 9).Value = 3122
$ws.Cells.Item(32, 11).Value = 3122
$ws.Cells.Item(32, 13).Value = -2835
$ws.Cells.Item(34, 8).Value = 32599.1
$ws.Cells.Item(34, 9).Value = 28332.5
$ws.Cells.Item(34, 11).Value = 28332.5
$ws.Cells.Item(34, 13).Value = -28061.5
$ws.Cells.Item(45, 8).Value = 4442.25
$ws.Cells.Item(45, 9).Value = 4441.5
$ws.Cells.Item(45, 10).Value = 4443
$ws.Cells.Item(45, 11).Value = 4441.5
$ws.Cells.Item(45, 12).Value = 4443
$ws.Cells.Item(45, 13).Value = -4064.5
$ws.Cells.Item(45, 14).Value = -5197
$ws.Cells.Item(61, 8).Value = 2135.7083
$ws.Cells.Item(61, 9).Value = 2135.7083
$ws.Cells.Item(61, 11).Value = 2135.7083
$ws.Cells.Item(61, 13).Value = -1923.7083
$ws.Cells.Item(74, 8).Value = 1562.0938
$ws.Cells.Item(74, 9).Value = 1373.6296
$ws.Cells.Item(74, 10).Value = 2579.8
$ws.Cells.Item(74, 11).Value = 1373.6296
$ws.Cells.Item(74, 12).Value = 2579.8
$ws.Cells.Item(74, 13).Value = -499.6296
$ws.Cells.Item(74, 14).Value = -4327.8
$ws.Cells.Item(77, 8).Value = 1562.0938
$ws.Cells.Item(77, 9).Value = 1373.6296
$ws.Cells.Item(77, 10).Value = 2579.8
$ws.Cells.Item(77, 11).Value = 6868.148
$ws.Cells.Item(77, 12).Value = 12899
$ws.Cells.Item(77, 13).Value = -2500.148
$ws.Cells.Item(77, 14).Value = -21635
$ws.Cells.Item(116, 8).Value = 2269
$ws.Cells.Item(116, 9).Value = 2674.75
$ws.Cells.Item(116, 10).Value = 1457.5
$ws.Cells.Item(116, 11).Value = 2674.75
$ws.Cells.Item(116, 12).Value = 1457.5
$ws.Cells.Item(116, 13).Value = -380.75
$ws.Cells.Item(116, 14).Value = -6045.5
$ws.Cells.Item(130, 8).Value = 12333.333
$ws.Cells.Item(130, 10).Value = 12333.333
$ws.Cells.Item(130, 12).Value = 12333.333
$ws.Cells.Item(130, 14).Value = -22373.333
$ws.Cells.Item(136, 8).Value = 2135.7083
$ws.Cells.Item(136, 9).Value = 2135.7083
$ws.Cells.Item(136, 11).Value = 6407.124899999999
$ws.Cells.Item(136, 13).Value = -3857.124899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2269
$ws.Cells.Item(3, 9).Value = 2674.75
$ws.Cells.Item(3, 10).Value = 1457.5
$ws.Cells.Item(3, 11).Value = 2674.75
$ws.Cells.Item(3, 12).Value = 1457.5
$ws.Cells.Item(3, 13).Value = -2560.75
$ws.Cells.Item(3, 14).Value = -1685.5
$ws.Cells.Item(4, 8).Value = 1017.1667
$ws.Cells.Item(4, 9).Value = 620.6
$ws.Cells.Item(4, 11).Value = 620.6
$ws.Cells.Item(4, 13).Value = -505.6
$ws.Cells.Item(22, 8).Value = 465.5
$ws.Cells.Item(22, 9).Value = 464
$ws.Cells.Item(22, 10).Value = 470
$ws.Cells.Item(22, 11).Value = 464
$ws.Cells.Item(22, 12).Value = 470
$ws.Cells.Item(22, 13).Value = -291
$ws.Cells.Item(22, 14).Value = -816
$ws.Cells.Item(49, 8).Value = 14000
$ws.Cells.Item(49, 9).Value = 14000
$ws.Cells.Item(49, 11).Value = 14000
$ws.Cells.Item(49, 13).Value = -13761
$ws.Cells.Item(74, 8).Value = 29760.75
$ws.Cells.Item(74, 10).Value = 19679.5
$ws.Cells.Item(74, 12).Value = 19679.5
$ws.Cells.Item(74, 14).Value = -21551.5
$ws.Cells.Item(77, 8).Value = 29760.75
$ws.Cells.Item(77, 10).Value = 19679.5
$ws.Cells.Item(77, 12).Value = 59038.5
$ws.Cells.Item(77, 14).Value = -68398.5
$ws.Cells.Item(86, 8).Value = 1703.5
$ws.Cells.Item(86, 9).Value = 1703.5
$ws.Cells.Item(86, 11).Value = 1703.5
$ws.Cells.Item(86, 13).Value = -580.5
$ws.Cells.Item(88, 8).Value = 4300
$ws.Cells.Item(88, 10).Value = 4300
$ws.Cells.Item(88, 12).Value = 4300
$ws.Cells.Item(88, 14).Value = -5112
$ws.Cells.Item(89, 8).Value = 1703.5
$ws.Cells.Item(89, 9).Value = 1703.5
$ws.Cells.Item(89, 11).Value = 8517.5
$ws.Cells.Item(89, 13).Value = -2901.5
$ws.Cells.Item(91, 8).Value = 4300
$ws.Cells.Item(91, 10).Value = 4300
$ws.Cells.Item(91, 12).Value = 4300
$ws.Cells.Item(91, 14).Value = -7108
$ws.Cells.Item(94, 8).Value = 2516.9565
$ws.Cells.Item(94, 9).Value = 1013.7692
$ws.Cells.Item(94, 10).Value = 4471.1
$ws.Cells.Item(94, 11).Value = 1013.7692
$ws.Cells.Item(94, 12).Value = 4471.1
$ws.Cells.Item(94, 13).Value = -562.7692
$ws.Cells.Item(94, 14).Value = -5373.1
$ws.Cells.Item(99, 8).Value = 4666.6665
$ws.Cells.Item(99, 9).Value = 4666.6665
$ws.Cells.Item(99, 11).Value = 4666.6665
$ws.Cells.Item(99, 13).Value = -3168.6665
$ws.Cells.Item(105, 8).Value = 2323.5
$ws.Cells.Item(105, 9).Value = 2245.5454
$ws.Cells.Item(105, 10).Value = 2389.4614
$ws.Cells.Item(105, 11).Value = 2245.5454
$ws.Cells.Item(105, 12).Value = 2389.4614
$ws.Cells.Item(105, 13).Value = -498.5454
$ws.Cells.Item(105, 14).Value = -5883.4614
$ws.Cells.Item(107, 8).Value = 867.6667
$ws.Cells.Item(107, 9).Value = 786.7857
$ws.Cells.Item(107, 11).Value = 786.7857
$ws.Cells.Item(107, 13).Value = 1133.2143

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 476
$ws.Cells.Item(16, 9).Value = 476
$ws.Cells.Item(16, 11).Value = 476
$ws.Cells.Item(16, 13).Value = -189
$ws.Cells.Item(31, 8).Value = 1907
$ws.Cells.Item(31, 9).Value = 1462.0667
$ws.Cells.Item(31, 10).Value = 2351.9333
$ws.Cells.Item(31, 11).Value = 1462.0667
$ws.Cells.Item(31, 12).Value = 2351.9333
$ws.Cells.Item(31, 13).Value = -1167.0667
$ws.Cells.Item(31, 14).Value = -2941.9333
$ws.Cells.Item(34, 8).Value = 1907
$ws.Cells.Item(34, 9).Value = 1462.0667
$ws.Cells.Item(34, 10).Value = 2351.9333
$ws.Cells.Item(34, 11).Value = 1462.0667
$ws.Cells.Item(34, 12).Value = 2351.9333
$ws.Cells.Item(34, 13).Value = -1260.0667
$ws.Cells.Item(34, 14).Value = -2755.9333
$ws.Cells.Item(35, 8).Value = 2808
$ws.Cells.Item(35, 9).Value = 2808
$ws.Cells.Item(35, 11).Value = 2808
$ws.Cells.Item(35, 13).Value = -2514
$ws.Cells.Item(36, 8).Value = 0
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 12).Value = $null
$ws.Cells.Item(36, 14).Value = 0
$ws.Cells.Item(38, 8).Value = 8000
$ws.Cells.Item(38, 9).Value = 8000
$ws.Cells.Item(38, 11).Value = 8000
$ws.Cells.Item(38, 13).Value = -7623
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 12).Value = $null
$ws.Cells.Item(40, 14).Value = 0
$ws.Cells.Item(46, 8).Value = 8000
$ws.Cells.Item(46, 9).Value = 8000
$ws.Cells.Item(46, 11).Value = 8000
$ws.Cells.Item(46, 13).Value = -7789
$ws.Cells.Item(54, 8).Value = 63500
$ws.Cells.Item(54, 10).Value = 63500
$ws.Cells.Item(54, 12).Value = 63500
$ws.Cells.Item(54, 14).Value = -64816
$ws.Cells.Item(58, 8).Value = 2122.1904
$ws.Cells.Item(58, 9).Value = 2175.6316
$ws.Cells.Item(58, 10).Value = 1614.5
$ws.Cells.Item(58, 11).Value = 2175.6316
$ws.Cells.Item(58, 12).Value = 1614.5
$ws.Cells.Item(58, 13).Value = -1972.6316
$ws.Cells.Item(58, 14).Value = -2020.5
$ws.Cells.Item(99, 8).Value = 170000
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 11).Value = 0
$ws.Cells.Item(99, 13).Value = $null
$ws.Cells.Item(105, 8).Value = 2581.2307
$ws.Cells.Item(105, 10).Value = 4339.2
$ws.Cells.Item(105, 12).Value = 4339.2
$ws.Cells.Item(105, 14).Value = -7833.2
$ws.Cells.Item(107, 8).Value = 1649.35
$ws.Cells.Item(107, 9).Value = 1460.5385
$ws.Cells.Item(107, 11).Value = 1460.5385
$ws.Cells.Item(107, 13).Value = 459.4614999999999
$ws.Cells.Item(113, 8).Value = 476
$ws.Cells.Item(113, 9).Value = 476
$ws.Cells.Item(113, 11).Value = 476
$ws.Cells.Item(113, 13).Value = 1694
$ws.Cells.Item(122, 8).Value = 1475
$ws.Cells.Item(122, 9).Value = 1475
$ws.Cells.Item(122, 11).Value = 4425
$ws.Cells.Item(122, 13).Value = -1975
$ws.Cells.Item(126, 8).Value = 170000
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 13).Value = $null
$ws.Cells.Item(134, 8).Value = 2797.9412
$ws.Cells.Item(134, 9).Value = 2785.3125
$ws.Cells.Item(134, 10).Value = 3000
$ws.Cells.Item(134, 11).Value = 8355.9375
$ws.Cells.Item(134, 12).Value = 9000
$ws.Cells.Item(134, 13).Value = -5820.9375
$ws.Cells.Item(134, 14).Value = -14070
$ws.Cells.Item(136, 8).Value = 2122.1904
$ws.Cells.Item(136, 9).Value = 2175.6316
$ws.Cells.Item(136, 10).Value = 1614.5
$ws.Cells.Item(136, 11).Value = 6526.8948
$ws.Cells.Item(136, 12).Value = 4843.5
$ws.Cells.Item(136, 13).Value = -3976.8948
$ws.Cells.Item(136, 14).Value = -9943.5
$ws.Cells.Item(141, 8).Value = 404372.8
$ws.Cells.Item(141, 10).Value = 22000
$ws.Cells.Item(141, 12).Value = 22000
$ws.Cells.Item(141, 14).Value = -32360

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 46435132
$ws.Cells.Item(4, 9).Value = 2962901.8
$ws.Cells.Item(4, 10).Value = 466666700
$ws.Cells.Item(4, 11).Value = 8888705.399999999
$ws.Cells.Item(4, 12).Value = 1400000100
$ws.Cells.Item(4, 13).Value = -8888593.399999999
$ws.Cells.Item(4, 14).Value = -1400000324
$ws.Cells.Item(5, 8).Value = 809.3333
$ws.Cells.Item(5, 9).Value = 747
$ws.Cells.Item(5, 10).Value = 996.3333
$ws.Cells.Item(5, 11).Value = 2241
$ws.Cells.Item(5, 12).Value = 2988.9999
$ws.Cells.Item(5, 13).Value = -2129
$ws.Cells.Item(5, 14).Value = -3212.9999
$ws.Cells.Item(9, 8).Value = 3278.6
$ws.Cells.Item(9, 10).Value = 2598.75
$ws.Cells.Item(9, 12).Value = 7796.25
$ws.Cells.Item(9, 14).Value = -8244.25
$ws.Cells.Item(12, 8).Value = 262.70834
$ws.Cells.Item(12, 10).Value = 191.15384
$ws.Cells.Item(12, 12).Value = 573.4615200000001
$ws.Cells.Item(12, 14).Value = -919.4615200000001
$ws.Cells.Item(32, 8).Value = 1000
$ws.Cells.Item(32, 9).Value = 1000
$ws.Cells.Item(32, 11).Value = 3000
$ws.Cells.Item(32, 13).Value = -2717
$ws.Cells.Item(51, 8).Value = 4643
$ws.Cells.Item(51, 9).Value = 312
$ws.Cells.Item(51, 11).Value = 936
$ws.Cells.Item(51, 13).Value = -476
$ws.Cells.Item(68, 8).Value = 2551.6667
$ws.Cells.Item(68, 9).Value = 1566.7142
$ws.Cells.Item(68, 11).Value = 4700.142599999999
$ws.Cells.Item(68, 13).Value = -3889.142599999999
$ws.Cells.Item(71, 8).Value = 2551.6667
$ws.Cells.Item(71, 9).Value = 1566.7142
$ws.Cells.Item(71, 11).Value = 14100.4278
$ws.Cells.Item(71, 13).Value = -10044.4278
$ws.Cells.Item(92, 8).Value = 231.25
$ws.Cells.Item(92, 9).Value = 283
$ws.Cells.Item(92, 10).Value = 223.85715
$ws.Cells.Item(92, 11).Value = 849
$ws.Cells.Item(92, 12).Value = 671.5714499999999
$ws.Cells.Item(92, 13).Value = 399
$ws.Cells.Item(92, 14).Value = -3167.57145
$ws.Cells.Item(93, 8).Value = 1316.6666
$ws.Cells.Item(93, 10).Value = 1400
$ws.Cells.Item(93, 12).Value = 4200
$ws.Cells.Item(93, 14).Value = -7944
$ws.Cells.Item(94, 8).Value = 6932
$ws.Cells.Item(94, 9).Value = 0
$ws.Cells.Item(94, 11).Value = 0
$ws.Cells.Item(94, 13).Value = $null
$ws.Cells.Item(97, 8).Value = 2043.3334
$ws.Cells.Item(97, 9).Value = 2355.7144
$ws.Cells.Item(97, 10).Value = 950
$ws.Cells.Item(97, 11).Value = 7067.1432
$ws.Cells.Item(97, 12).Value = 2850
$ws.Cells.Item(97, 13).Value = -6571.1432
$ws.Cells.Item(97, 14).Value = -3842
$ws.Cells.Item(98, 8).Value = 1000
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 14).Value = $null
$ws.Cells.Item(99, 8).Value = 0
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 11).Value = 0
$ws.Cells.Item(99, 13).Value = $null
$ws.Cells.Item(100, 8).Value = 2000
$ws.Cells.Item(100, 10).Value = 2000
$ws.Cells.Item(100, 12).Value = 6000
$ws.Cells.Item(100, 14).Value = -7622
$ws.Cells.Item(103, 8).Value = 255.77777
$ws.Cells.Item(103, 9).Value = 94.333336
$ws.Cells.Item(103, 10).Value = 336.5
$ws.Cells.Item(103, 11).Value = 283.000008
$ws.Cells.Item(103, 12).Value = 1009.5
$ws.Cells.Item(103, 13).Value = 595.999992
$ws.Cells.Item(103, 14).Value = -2767.5
$ws.Cells.Item(134, 8).Value = 571.2222
$ws.Cells.Item(134, 9).Value = 571.2222
$ws.Cells.Item(134, 11).Value = 1713.6666
$ws.Cells.Item(134, 13).Value = 3356.3334
$ws.Cells.Item(135, 8).Value = 809.3333
$ws.Cells.Item(135, 9).Value = 747
$ws.Cells.Item(135, 10).Value = 996.3333
$ws.Cells.Item(135, 11).Value = 6723
$ws.Cells.Item(135, 12).Value = 8966.9997
$ws.Cells.Item(135, 13).Value = -4188
$ws.Cells.Item(135, 14).Value = -14036.9997
$ws.Cells.Item(140, 8).Value = 10824.6
$ws.Cells.Item(140, 9).Value = 4513
$ws.Cells.Item(140, 10).Value = 15032.333
$ws.Cells.Item(140, 11).Value = 13539
$ws.Cells.Item(140, 12).Value = 45096.999
$ws.Cells.Item(140, 13).Value = -8359
$ws.Cells.Item(140, 14).Value = -55456.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 7740.067
$ws.Cells.Item(70, 10).Value = 8845
$ws.Cells.Item(70, 12).Value = 8845
$ws.Cells.Item(70, 14).Value = -9385
$ws.Cells.Item(73, 8).Value = 7740.067
$ws.Cells.Item(73, 10).Value = 8845
$ws.Cells.Item(73, 12).Value = 8845
$ws.Cells.Item(73, 14).Value = -10717
$ws.Cells.Item(80, 8).Value = 22041.5
$ws.Cells.Item(80, 9).Value = 6166.5
$ws.Cells.Item(80, 10).Value = 27333.166
$ws.Cells.Item(80, 11).Value = 6166.5
$ws.Cells.Item(80, 12).Value = 27333.166
$ws.Cells.Item(80, 13).Value = -5168.5
$ws.Cells.Item(80, 14).Value = -29329.166
$ws.Cells.Item(83, 8).Value = 22041.5
$ws.Cells.Item(83, 9).Value = 6166.5
$ws.Cells.Item(83, 10).Value = 27333.166
$ws.Cells.Item(83, 11).Value = 30832.5
$ws.Cells.Item(83, 12).Value = 136665.83
$ws.Cells.Item(83, 13).Value = -25840.5
$ws.Cells.Item(83, 14).Value = -146649.83
$ws.Cells.Item(102, 8).Value = 3304
$ws.Cells.Item(102, 9).Value = 2510.5
$ws.Cells.Item(102, 10).Value = 4210.857
$ws.Cells.Item(102, 11).Value = 2510.5
$ws.Cells.Item(102, 12).Value = 4210.857
$ws.Cells.Item(102, 13).Value = -888.5
$ws.Cells.Item(102, 14).Value = -7454.857
$ws.Cells.Item(122, 8).Value = 2249.0715
$ws.Cells.Item(122, 9).Value = 1770.7778
$ws.Cells.Item(122, 10).Value = 3110
$ws.Cells.Item(122, 11).Value = 5312.3334
$ws.Cells.Item(122, 12).Value = 9330
$ws.Cells.Item(122, 13).Value = -2862.3334
$ws.Cells.Item(122, 14).Value = -14230
$ws.Cells.Item(132, 8).Value = 2650.8572
$ws.Cells.Item(132, 9).Value = 2650.8572
$ws.Cells.Item(132, 11).Value = 7952.571599999999
$ws.Cells.Item(132, 13).Value = -5422.571599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(43, 8).Value = 3792080.8
$ws.Cells.Item(43, 10).Value = 5683121
$ws.Cells.Item(43, 12).Value = 5683121
$ws.Cells.Item(43, 14).Value = -5683507
$ws.Cells.Item(55, 8).Value = 2740.9092
$ws.Cells.Item(55, 9).Value = 483.22223
$ws.Cells.Item(55, 11).Value = 483.22223
$ws.Cells.Item(55, 13).Value = -310.22223
$ws.Cells.Item(68, 8).Value = 2405.75
$ws.Cells.Item(68, 9).Value = 2405.75
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 11).Value = 2405.75
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 13).Value = $null
$ws.Cells.Item(68, 14).Value = -1656.75
$ws.Cells.Item(71, 8).Value = 2405.75
$ws.Cells.Item(71, 9).Value = 2405.75
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 12028.75
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 13).Value = $null
$ws.Cells.Item(71, 14).Value = -8284.75
$ws.Cells.Item(82, 8).Value = 2663.1667
$ws.Cells.Item(82, 9).Value = 2541.6365
$ws.Cells.Item(82, 10).Value = 4000
$ws.Cells.Item(82, 11).Value = 2541.6365
$ws.Cells.Item(82, 12).Value = 4000
$ws.Cells.Item(82, 13).Value = -2180.6365
$ws.Cells.Item(82, 14).Value = -4722
$ws.Cells.Item(85, 8).Value = 2663.1667
$ws.Cells.Item(85, 9).Value = 2541.6365
$ws.Cells.Item(85, 10).Value = 4000
$ws.Cells.Item(85, 11).Value = 2541.6365
$ws.Cells.Item(85, 12).Value = 4000
$ws.Cells.Item(85, 13).Value = -1293.6365
$ws.Cells.Item(85, 14).Value = -6496
$ws.Cells.Item(100, 8).Value = 2724.875
$ws.Cells.Item(100, 9).Value = 1933
$ws.Cells.Item(100, 10).Value = 3200
$ws.Cells.Item(100, 11).Value = 1933
$ws.Cells.Item(100, 12).Value = 3200
$ws.Cells.Item(100, 13).Value = -1392
$ws.Cells.Item(100, 14).Value = -4282
$ws.Cells.Item(122, 8).Value = 3823.6667
$ws.Cells.Item(122, 9).Value = 4321.5713
$ws.Cells.Item(122, 11).Value = 12964.7139
$ws.Cells.Item(122, 13).Value = -10514.7139
$ws.Cells.Item(132, 8).Value = 3184.2307
$ws.Cells.Item(132, 9).Value = 2742.3809
$ws.Cells.Item(132, 10).Value = 5040
$ws.Cells.Item(132, 11).Value = 8227.1427
$ws.Cells.Item(132, 12).Value = 15120
$ws.Cells.Item(132, 13).Value = -5697.1427
$ws.Cells.Item(132, 14).Value = -20180
$ws.Cells.Item(136, 8).Value = 1805.8948
$ws.Cells.Item(136, 9).Value = 1375.75
$ws.Cells.Item(136, 11).Value = 4127.25
$ws.Cells.Item(136, 13).Value = -1577.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(25, 8).Value = 8998.5
$ws.Cells.Item(25, 10).Value = 8998.5
$ws.Cells.Item(25, 12).Value = 8998.5
$ws.Cells.Item(25, 14).Value = -9584.5
$ws.Cells.Item(34, 8).Value = 35249.668
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 14).Value = $null
$ws.Cells.Item(40, 8).Value = 49999
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 14).Value = $null
$ws.Cells.Item(42, 8).Value = 39997.5
$ws.Cells.Item(42, 9).Value = 39997.5
$ws.Cells.Item(42, 11).Value = 39997.5
$ws.Cells.Item(42, 13).Value = -39619.5
$ws.Cells.Item(81, 8).Value = 3158.4167
$ws.Cells.Item(81, 9).Value = 1780.1
$ws.Cells.Item(81, 10).Value = 10050
$ws.Cells.Item(81, 11).Value = 3560.2
$ws.Cells.Item(81, 12).Value = 20100
$ws.Cells.Item(81, 13).Value = -2499.2
$ws.Cells.Item(81, 14).Value = -22222
$ws.Cells.Item(84, 8).Value = 3158.4167
$ws.Cells.Item(84, 9).Value = 1780.1
$ws.Cells.Item(84, 10).Value = 10050
$ws.Cells.Item(84, 11).Value = 17801
$ws.Cells.Item(84, 12).Value = 100500
$ws.Cells.Item(84, 13).Value = -12497
$ws.Cells.Item(84, 14).Value = -111108
$ws.Cells.Item(107, 8).Value = 709.62964
$ws.Cells.Item(107, 9).Value = 495.53845
$ws.Cells.Item(107, 10).Value = 908.4286
$ws.Cells.Item(107, 11).Value = 1486.61535
$ws.Cells.Item(107, 12).Value = 2725.2858
$ws.Cells.Item(107, 13).Value = 433.38465
$ws.Cells.Item(107, 14).Value = -6565.2858
$ws.Cells.Item(113, 8).Value = 509
$ws.Cells.Item(113, 9).Value = 421
$ws.Cells.Item(113, 10).Value = 575
$ws.Cells.Item(113, 11).Value = 1263
$ws.Cells.Item(113, 12).Value = 1725
$ws.Cells.Item(113, 13).Value = 907
$ws.Cells.Item(113, 14).Value = -6065
$ws.Cells.Item(122, 8).Value = 5859.242
$ws.Cells.Item(122, 9).Value = 6814.7
$ws.Cells.Item(122, 11).Value = 20444.1
$ws.Cells.Item(122, 13).Value = -17994.1
$ws.Cells.Item(126, 8).Value = 1431.8572
$ws.Cells.Item(126, 9).Value = 1448.5
$ws.Cells.Item(126, 10).Value = 1429.0834
$ws.Cells.Item(126, 11).Value = 4345.5
$ws.Cells.Item(126, 12).Value = 4287.2502
$ws.Cells.Item(126, 13).Value = -1875.5
$ws.Cells.Item(126, 14).Value = -9227.2502
$ws.Cells.Item(132, 8).Value = 1145.8182
$ws.Cells.Item(132, 9).Value = 1180.381
$ws.Cells.Item(132, 10).Value = 420
$ws.Cells.Item(132, 11).Value = 3541.143
$ws.Cells.Item(132, 12).Value = 1260
$ws.Cells.Item(132, 13).Value = -1011.143
$ws.Cells.Item(132, 14).Value = -6320
$ws.Cells.Item(141, 8).Value = 67499.164
$ws.Cells.Item(141, 10).Value = 67499.164
$ws.Cells.Item(141, 12).Value = 67499.164
$ws.Cells.Item(141, 14).Value = -77859.164
